# Convert the "First Year" (C) and "Last Year" (D) columns from full date
# serials (formatted as yyyy-mm-dd / custom datetime formats) into plain
# 4-digit year numbers, and strip the date number-formatting from those
# cells so they display as plain numbers again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel-correct "serial date -> calendar year" conversion. .NET's
# [datetime]::FromOADate() uses the same OLE Automation date epoch as
# Excel, but it does not reproduce Excel's (incorrect) belief that 1900
# was a leap year, so serials 1-59 land one day early. Nudge those by a
# day so the year matches what Excel itself displays.
function Get-ExcelYear($serial) {
    $s = [double]$serial
    if ($s -ge 1 -and $s -lt 60) {
        $s = $s + 1
    }
    return ([datetime]::FromOADate($s)).Year
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 146 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cSerial = $cCell.Value2
    $dSerial = $dCell.Value2

    if ($cSerial -eq $null -and $dSerial -eq $null) {
        continue
    }

    $cYear = Get-ExcelYear $cSerial
    $dYear = Get-ExcelYear $dSerial

    $rangeAddr = "C$r" + ":" + "D$r"
    $ws.Range($rangeAddr).ClearFormats()

    $cCell.Value = $cYear
    $dCell.Value = $dYear
}
